$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1889129372953526,0.05145540628087652,0.510633136522584,1.150850776514954,0.244051625780531,0.007143758385484129,0.2765385556732914,0.2091803144137915,0.2790927646119826),
    @(0.1147728127515536,0.1276683352405244,0.5556586317539761,1.153506061810903,0.2280370556961183,0.007145814755611407,0.2702865095000494,0.2118287548144138,0.2706881858766149),
    @(0.08605073027678881,0.1840294754951445,0.5869376741634227,1.150640452582486,0.2152958185488645,0.007147202478436689,0.2660713951214769,0.2088751067174658,0.2637329729143062),
    @(0.088634560928939,0.2083061261039598,0.5977682494760557,1.143075964402154,0.2085478735509382,0.007148026877511104,0.2635639320690886,0.2011613001046918,0.2591364487608195),
    @(0.09208648359710619,0.2167500571965187,0.6003016410720982,1.137487374840272,0.2056058570367308,0.007148451054963098,0.2622743362376215,0.195468654388216,0.2567448189688482),
    @(0.09208648359710619,0.2167500571965187,0.6003016410720982,1.137487374840272,0.2056058570367308,0.007148451054963098,0.2622743362376215,0.195468654388216,0.2567448189688482),
    @(0.09022229300813775,0.2123535318682941,0.5973365438421244,1.136442058871971,0.2063307752735777,0.007148436850666838,0.2623200524810719,0.194407825461798,0.256885550488286),
    @(0.09051380346537016,0.1619028615249483,0.5675658088136332,1.134805100903795,0.2168173355540869,0.007147591079763722,0.2648913898532926,0.1927765353954863,0.261622082692071),
    @(0.1921047816904505,0.04659520757678177,0.5032842372182631,1.139980812592168,0.2428591421129651,0.007144187940705799,0.2752374863347197,0.1981256957906939,0.2765843321225072),
    @(0.3197146668046357,0.07109551978376344,0.4410736289854761,1.143570937740073,0.2678713340348098,0.007139662611256327,0.289001806760822,0.2018601469471587,0.2927330823388139),
    @(0.3803321028281249,0.1217034295276755,0.4084484459954485,1.141918122458789,0.279794686496178,0.007137154114416127,0.2958388203751171,0.200220603374042,0.3005628947597386),
    @(0.3975149561865657,0.1303587383554133,0.3827277756885337,1.140090616236944,0.2890403932338443,0.00713545529950798,0.2980199358906649,0.1983925647331067,0.3064555396357527),
    @(0.3779150855298699,0.1087606416060044,0.3760614111285185,1.137498027912709,0.2906458897652517,0.007135465999384094,0.2987805122313302,0.1957622157117308,0.3068903188376896),
    @(0.3567680278826698,0.08937561498686801,0.3784122077525797,1.133661149489371,0.2883002572875591,0.007136153901236136,0.2979917845436875,0.1918557990537162,0.3045460271897957),
    @(0.3488996029039972,0.08128804930620263,0.3819807263227016,1.137419184907291,0.2883224769554886,0.007136151406264199,0.2985906216765174,0.1956738137861299,0.3050762539993411),
    @(0.3417990623091169,0.07387309515439706,0.3863174116912869,1.14229769903722,0.2884660860991557,0.007136095401005098,0.2992396713791387,0.2006309080445968,0.3058744442228326),
    @(0.3049870983268566,0.04376019637099582,0.4054998074787127,1.14185053555035,0.2808944241919179,0.007137660188404726,0.2942614783753689,0.2001522191454512,0.3006139270167011),
    @(0.2887170202833134,0.03210620233011515,0.4197886867432128,1.144048218325694,0.2762118318930486,0.007138499047076671,0.2904561912760948,0.2023705197569689,0.297864715120648),
    @(0.3274313816481048,0.07090671029021982,0.4192655736466873,1.143618564822357,0.2762569552284753,0.007138244421478918,0.2896706017494096,0.2019352018892712,0.2980652247025294),
    @(0.425825288677754,0.1597547333116663,0.3898230719570778,1.146045718935838,0.2885275398514514,0.007135025504571312,0.3026174511751653,0.204443681567211,0.3074292079699993),
    @(0.4884954286003906,0.2120515798159485,0.3643775833863607,1.151577979722563,0.3004793886203635,0.007131873839731551,0.3126558654020301,0.2101026045794573,0.3171765506049416),
    @(0.5328512648782513,0.2482655875382686,0.3466216043032792,1.157559256286586,0.3096608336066091,0.007129319943871478,0.3204277597197207,0.2162065739262544,0.325138868297276),
    @(0.5474933645431455,0.2598325462924789,0.3419480036746411,1.161776743503445,0.3130815587770856,0.007128339617714072,0.3234149321151464,0.2205008683162254,0.3284311882861853),
    @(0.4473040772317025,0.1751665974230721,0.3907641001539978,1.16494078617554,0.2952269658667213,0.007133214733652158,0.3085927413845636,0.2236643829870655,0.3151246870219502),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}